$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 708
$ws1.Range("F3").Value = 13569
$ws1.Range("F4").Value = 13809
$ws1.Range("F7").Value = 1366
$ws1.Range("F8").Value = 5749
$ws1.Range("F9").Value = 957
$ws1.Range("F10").Value = 559
$ws1.Range("F14").Value = 1499
$ws1.Range("F15").Value = 410
$ws1.Range("F16").Value = 2122
$ws1.Range("F17").Value = 1147
$ws1.Range("F18").Value = 1744
$ws1.Range("F19").Value = 902
$ws1.Range("F21").Value = 2232
$ws1.Range("F22").Value = 543
$ws1.Range("F23").Value = 777
$ws1.Range("F24").Value = 3204
$ws1.Range("F26").Value = 288
$ws1.Range("F27").Value = 2267
$ws1.Range("F28").Value = 58
$ws1.Range("F31").Value = 1744
$ws1.Range("F33").Value = 1290
$ws1.Range("F35").Value = 125
$ws1.Range("F36").Value = 4548
$ws1.Range("F37").Value = 4656
$ws1.Range("F39").Value = 148
$ws1.Range("F40").Value = 645
$ws1.Range("F42").Value = 3232
$ws1.Range("F45").Value = 320
$ws1.Range("F47").Value = 57
$ws1.Range("F48").Value = 4380
$ws1.Range("F49").Value = 256

# Sheet 3: 本地生活 (index 3)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 7319
$ws3.Range("F4").Value = 609

# Sheet 4: 全部类型 (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 708
$ws4.Range("F4").Value = 609
$ws4.Range("F6").Value = 13569
$ws4.Range("F7").Value = 13809
$ws4.Range("F10").Value = 5749
$ws4.Range("F11").Value = 957
$ws4.Range("F16").Value = 1499
$ws4.Range("F17").Value = 410
$ws4.Range("F18").Value = 2122
$ws4.Range("F19").Value = 1147
$ws4.Range("F20").Value = 1744
$ws4.Range("F21").Value = 902
$ws4.Range("F22").Value = 543
$ws4.Range("F23").Value = 3204
$ws4.Range("F24").Value = 288
$ws4.Range("F25").Value = 58
$ws4.Range("F28").Value = 1744
$ws4.Range("F31").Value = 1290
$ws4.Range("F34").Value = 4548
$ws4.Range("F35").Value = 4656
$ws4.Range("F38").Value = 148
$ws4.Range("F39").Value = 645
$ws4.Range("F41").Value = 3232
$ws4.Range("F44").Value = 320
$ws4.Range("F46").Value = 57
$ws4.Range("F47").Value = 4380
$ws4.Range("F48").Value = 256
